$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.876.43"
$ws.Range("E2").Value = "'  -3.14%  "
$ws.Range("D3").Value = "'3.288.73"
$ws.Range("E3").Value = "'  -5.26%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'227.17"
$ws.Range("E5").Value = "'  -5.43%  "
$ws.Range("E6").Value = "'  -5.58%  "
$ws.Range("E7").Value = "'  -8.05%  "
$ws.Range("D8").Value = "'0.376"
$ws.Range("E8").Value = "'  -7.10%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "'  +0.03%  "
$ws.Range("D10").Value = "'0.929"
$ws.Range("E10").Value = "'  -8.12%  "
$ws.Range("D11").Value = "'3.284.18"
$ws.Range("E11").Value = "'  -5.35%  "
$ws.Range("D12").Value = "'41.29"
$ws.Range("E12").Value = "'  -2.17%  "
$ws.Range("E13").Value = "'  -3.86%  "
$ws.Range("D14").Value = "'5.90"
$ws.Range("E14").Value = "'  -3.91%  "
$ws.Range("D15").Value = "'91.733.37"
$ws.Range("E15").Value = "'  -3.10%  "
$ws.Range("D16").Value = "'3.898.27"
$ws.Range("E16").Value = "'  -5.50%  "
$ws.Range("D17").Value = "'0.0000239"
$ws.Range("E17").Value = "'  -6.65%  "
$ws.Range("D18").Value = "'7.96"
$ws.Range("E18").Value = "'  -6.65%  "
$ws.Range("D19").Value = "'3.284.81"
$ws.Range("E19").Value = "'  -5.39%  "
$ws.Range("D20").Value = "'16.94"
$ws.Range("E20").Value = "'  -5.56%  "
$ws.Range("D21").Value = "'10.55"
$ws.Range("E21").Value = "'  -7.83%  "
$ws.Range("D22").Value = "'3.36"
$ws.Range("E22").Value = "'  +5.31%  "
$ws.Range("D23").Value = "'482.34"
$ws.Range("E23").Value = "'  -4.19%  "
$ws.Range("D24").Value = "'0.433"
$ws.Range("E24").Value = "'  -16.80%  "
$ws.Range("D25").Value = "'0.0000176"
$ws.Range("E25").Value = "'  -8.84%  "
$ws.Range("D26").Value = "'6.08"
$ws.Range("E26").Value = "'  -8.41%  "
$ws.Range("D27").Value = "'88.62"
$ws.Range("E27").Value = "'  -2.78%  "
$ws.Range("D28").Value = "'11.58"
$ws.Range("E28").Value = "'  -4.59%  "
$ws.Range("D29").Value = "'3.463.92"
$ws.Range("E29").Value = "'  -5.30%  "
$ws.Range("E30").Value = "'  +0.04%  "
$ws.Range("D31").Value = "'10.80"
$ws.Range("E31").Value = "'  -8.79%  "
$ws.Range("D32").Value = "'0.136"
$ws.Range("E32").Value = "'  -1.98%  "
$ws.Range("E33").Value = "'  -6.70%  "
$ws.Range("E34").Value = "'  -0.03%  "
$ws.Range("E35").Value = "'  -7.89%  "
$ws.Range("D36").Value = "'27.62"
$ws.Range("E36").Value = "'  -10.57%  "
$ws.Range("D37").Value = "'0.512"
$ws.Range("E37").Value = "'  -9.96%  "
$ws.Range("D38").Value = "'535.68"
$ws.Range("E38").Value = "'  -0.75%  "
$ws.Range("E39").Value = "'  +0.02%  "
$ws.Range("E40").Value = "'  -7.31%  "
$ws.Range("E41").Value = "'  -3.46%  "
$ws.Range("E42").Value = "'  -7.81%  "
$ws.Range("D43").Value = "'0.846"
$ws.Range("E43").Value = "'  -9.48%  "
$ws.Range("D44").Value = "'23.72"
$ws.Range("E44").Value = "'  -1.47%  "
$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "'3.57"
$ws.Range("E45").Value = "'  +1.54%  "
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "'1.64"
$ws.Range("E46").Value = "'  -3.50%  "
$ws.Range("E47").Value = "'  -4.00%  "
$ws.Range("D48").Value = "'5.23"
$ws.Range("E48").Value = "'  -8.58%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'51.34"
$ws.Range("E49").Value = "'  -3.79%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.04"
$ws.Range("E50").Value = "'  -5.43%  "
$ws.Range("D51").Value = "'7.82"
$ws.Range("E51").Value = "'  -3.31%  "
